# Daily scrape update - 2025-12-27 03:32:23 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# (subtract the 5/6-char gridline padding Excel adds on top of ColumnWidth
# so the persisted <col width=".."/> lands on the exact target value)
$ws.Columns.Item(3).ColumnWidth = 42.166666666666664
$ws.Columns.Item(4).ColumnWidth = 55.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 55.166666666666664

# --- Row 2: update existing opportunity with freshly scraped values ---
$ws.Range("A2").Formula = "'1330867"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330867"
$ws.Range("C2").Value = "Social Media Director"
$ws.Range("D2").Value = "Başakşehir, Başak, 34490 Başakşehir/İstanbul, Türkiye"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Esen Isıtma Soğutma Elektrik İnşaat Sanayi ve Ticaret"

# --- Row 3: new opportunity ---
$ws.Range("A3").Formula = "'1330859"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1330859"
$ws.Range("C3").Value = "Export Sales Specialist"
$ws.Range("D3").Value = "Başakşehir, Başak, 34490 Başakşehir/İstanbul, Türkiye"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "2 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Esen Isıtma Soğutma Elektrik İnşaat Sanayi ve Ticaret"

# --- Row 4: new opportunity ---
$ws.Range("A4").Formula = "'1330856"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1330856"
$ws.Range("C4").Value = "Business Development Intern  (Long Term)"
$ws.Range("D4").Value = "Lahore, Pakistan"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Chughtaiz"

# --- Row 5: new opportunity ---
$ws.Range("A5").Formula = "'1330682"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1330682"
$ws.Range("C5").Value = "Web Master"
$ws.Range("D5").Value = "Başakşehir, Başak, 34490 Başakşehir/İstanbul, Türkiye"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Esen Isıtma Soğutma Elektrik İnşaat Sanayi ve Ticaret"

# --- Row 6: new opportunity ---
$ws.Range("A6").Formula = "'1322343"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1322343"
$ws.Range("C6").Value = "Electronics and Electrical Intern"
$ws.Range("D6").Value = "Manipal, Karnataka, India"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "21 applicants"
$ws.Range("G6").Value = "3 - 6 Months"
$ws.Range("H6").Value = "M.A.H.E."

Write-Output "applied daily scrape update"
